# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 293
$wsExpo.Range("F5").Value = 159
$wsExpo.Range("F6").Value = 127
$wsExpo.Range("F7").Value = 298
$wsExpo.Range("F8").Value = 218
$wsExpo.Range("F9").Value = 2057
$wsExpo.Range("F11").Value = 4921
$wsExpo.Range("F12").Value = 97
$wsExpo.Range("F13").Value = 342

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G2").Value = "不可售"

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("G2").Value = "不可售"
$wsAll.Range("F6").Value = 293
$wsAll.Range("F7").Value = 159
$wsAll.Range("F8").Value = 127
$wsAll.Range("F9").Value = 298
$wsAll.Range("F10").Value = 218
$wsAll.Range("F13").Value = 2057
$wsAll.Range("F15").Value = 4921
$wsAll.Range("F16").Value = 97
$wsAll.Range("F17").Value = 342
